$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete the experiment for p=2 (2^32 section): fill in the missing
# "100000" column (E) values for rows 24-27.
$ws.Range("E24").Value = 2633100
$ws.Range("E25").Value = 598649
$ws.Range("E26").Value = 8314
$ws.Range("E27").Value = 666

# Match style/formatting of the adjacent cells in the same rows (column D).
$ws.Range("D24:D27").Copy()
$ws.Range("E24:E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view state to match the saved selection/scroll position.
$ws.Range("F24").Select()
$excel.ActiveWindow.ScrollRow = 11
